$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the four data columns (A:D) to fit the new "best models" notebook content.
# Target OOXML <col width="..."> values (per the commit diff):
#   A -> 39.85, B -> 43.53, C -> 66.23, D -> 45.59
# Excel's ColumnWidth (character units) is fed through COM; the host snaps the
# stored width to the nearest displayable pixel, so we pick the ColumnWidth
# value that lands closest to each target after that rounding.
$ws.Columns.Item(1).ColumnWidth = 39
$ws.Columns.Item(2).ColumnWidth = 42.6666666666667
$ws.Columns.Item(3).ColumnWidth = 65.3333333333333
$ws.Columns.Item(4).ColumnWidth = 44.8333333333333
